$wb = $excel.ActiveWorkbook

# File identifiers for the new handback record
$newGuid     = "50f85b0a-7ac5-4813-b680-7f651db9ff13"
$newFileMd   = "$newGuid.md"
$newPathMd   = "e2e\$newGuid.md"
$newZhXlf    = "$newGuid.327286444e612566c607ce3177e91506ef919a4b.zh-cn.xlf"
$newDeXlf    = "$newGuid.327286444e612566c607ce3177e91506ef919a4b.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"

$dtOverview   = "2016-09-07 13:03:37"
$dtZhGenerate = "2016-09-07 13:03:24"
$dtZhHandback = "2016-09-07 13:04:33"
$dtDeGenerate = "2016-09-07 13:03:37"
$dtDeHandback = "2016-09-07 13:04:50"

$srcRepoUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fd7472234d33f3555c67f5f7749a3d57bf25f2e/e2e/$newFileMd"
$zhRepoUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/d9eba8db2a0ddc2439e807a6cdc886d46b41a008/e2e/$newFileMd"
$deRepoUrl   = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8e9f7cecdec3b8d21fedc2edbd54ed639a1a1915/e2e/$newFileMd"

# ---------------------------------------------------------------------------
# Sheet "Overview": add row 3
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFileMd
$wsOverview.Range("B3").Value = $newPathMd
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = $statusInSync
$wsOverview.Range("F3").Value = $statusInSync
$wsOverview.Range("G3").Value = $dtOverview

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $srcRepoUrl, "", "", $newPathMd) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": add row 3
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A3").Value = $newFileMd
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $statusInSync
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'True"
$wsZhCn.Range("G3").Value = $newZhXlf
$wsZhCn.Range("H3").Value = $dtZhGenerate
$wsZhCn.Range("I3").Value = $newFileMd
$wsZhCn.Range("J3").Value = $newZhXlf
$wsZhCn.Range("K3").Value = $dtZhHandback
$wsZhCn.Range("L3").Value = "'"
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("N3").Value = "'"
$wsZhCn.Range("O3").Value = "'False"
$wsZhCn.Range("P3").Value = "'"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $srcRepoUrl, "", "", $newFileMd) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $zhRepoUrl, "", "", $newFileMd) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": add row 3
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A3").Value = $newFileMd
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $statusInSync
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'True"
$wsDeDe.Range("G3").Value = $newDeXlf
$wsDeDe.Range("H3").Value = $dtDeGenerate
$wsDeDe.Range("I3").Value = $newFileMd
$wsDeDe.Range("J3").Value = $newDeXlf
$wsDeDe.Range("K3").Value = $dtDeHandback
$wsDeDe.Range("L3").Value = "'"
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("N3").Value = "'"
$wsDeDe.Range("O3").Value = "'False"
$wsDeDe.Range("P3").Value = "'"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $srcRepoUrl, "", "", $newFileMd) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $deRepoUrl, "", "", $newFileMd) | Out-Null
